$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 20:05"

$ws.Range("B4").Value = 1578245
$ws.Range("C4").Value = 7662
$ws.Range("D4").Value = 364683
$ws.Range("E4").Value = 1119628
$ws.Range("G4").Value = 401
$ws.Range("H4").Value = 93934

$ws.Range("B6").Value = 279524
$ws.Range("C6").Value = 721
$ws.Range("E6").Value = 54678
$ws.Range("G6").Value = 110
$ws.Range("H6").Value = 27888

$ws.Range("G8").Value = 363
$ws.Range("H8").Value = 35704

$ws.Range("B17").Value = 80081
$ws.Range("C17").Value = 969
$ws.Range("D17").Value = 40670
$ws.Range("E17").Value = 33384
$ws.Range("G17").Value = 115
$ws.Range("H17").Value = 6027

$ws.Range("B32").Value = 26004
$ws.Range("C32").Value = 941
$ws.Range("D32").Value = 11809
$ws.Range("E32").Value = 13962
$ws.Range("G32").Value = 6
$ws.Range("H32").Value = 233

$ws.Range("B33").Value = 24315
$ws.Range("C33").Value = 64
$ws.Range("E33").Value = 3274
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 1571

$ws.Range("B53").Value = 8281
$ws.Range("C53").Value = 14
$ws.Range("E53").Value = 8015

$ws.Range("B68").Value = 3971
$ws.Range("C68").Value = 13
$ws.Range("D68").Value = 3728
$ws.Range("E68").Value = 134

$ws.Range("B75").Value = 2939
$ws.Range("C75").Value = 84
$ws.Range("D75").Value = 2372
$ws.Range("E75").Value = 554

$ws.Range("B156").Value = 184
$ws.Range("C156").Value = 17
$ws.Range("E156").Value = 149
$ws.Range("G156").Value = 2
$ws.Range("H156").Value = 30
